$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scraper re-appended the same three match rows again (rows 3, 2, 4
# repeated, in that order) as new rows 5-7.
$newRows = @(
    @(" Abu Dhabi", " November 01 2020", "Super Kings won by 9 wickets (with 7 balls remaining)", "Kings XI Punjab", "Chennai Super Kings", "James Neesham ", "2", "3", "0", "0", "66.66"),
    @(" Abu Dhabi", " October 01 2020", "Mumbai won by 48 runs", "Kings XI Punjab", "Mumbai Indians", "James Neesham ", "7", "7", "0", "0", "100.00"),
    @(" Dubai (DSC)", " October 20 2020", "Kings XI won by 5 wickets (with 6 balls remaining)", "Kings XI Punjab", "Delhi Capitals", "James Neesham ", "10", "8", "0", "1", "125.00")
)

$startRow = 5
$lastRow = $startRow + $newRows.Length - 1

# The numeric-looking columns (G:K) must be stored as text, matching the
# rest of the sheet (ignoredErrors/numberStoredAsText covers A1:K7 too).
$ws.Range("G" + $startRow + ":K" + $lastRow).NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
